$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Rows 1-5: simple value replacements ---
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"
$t.Cell(4,1).Range.Text = "282"
$t.Cell(5,1).Range.Text = "0.00001"

# Row 6 (0.00059) is unchanged.

# --- Delete the two rows that held 0.00008 and 0.00005 (both were row 7) ---
$t.Rows.Item(7).Delete()
$t.Rows.Item(7).Delete()

# After the two deletions:
#   row 7 = 0.00010 (unchanged)
#   row 8 = 0.00010 -> 0.00004
#   row 9 = 0.00012 -> 0.00016
#   row 10 = 0.01223 -> 0.00017
$t.Cell(8,1).Range.Text = "0.00004"
$t.Cell(9,1).Range.Text = "0.00016"
$t.Cell(10,1).Range.Text = "0.00017"

# --- Insert two new rows after row 10 (before the old row 11 = 100.0) ---
$newRow1 = $t.Rows.Add($t.Rows.Item(11))
$t.Cell(11,1).Range.Text = "0.00020"
$newRow2 = $t.Rows.Add($t.Rows.Item(12))
$t.Cell(12,1).Range.Text = "0.03358"

# Rows 13 (100.0) through the row before the final three are unchanged.

# --- Final three rows: collapse the tab-separated runs into a single value ---
$lastCount = $t.Rows.Count
$t.Cell($lastCount - 2, 1).Range.Text = "99.89"
$t.Cell($lastCount - 1, 1).Range.Text = "0.03"
$t.Cell($lastCount, 1).Range.Text = "30"
